# Nuevo formato 15 jun 2021
#
# Updates the "Rescatables" roster (one student row swapped for another,
# one dropped-to-zero row removed) and refreshes the per-group pass/fail
# statistics on row 6 ("4ALCV") of the three "Estadisticos" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rescatables sheet: roster changes
#    - remove GARCIA / LINARES / ANDRES (NC 19330051920278)
#    - remove MORALES / CHIPAHUA / KARLA MICHEL (NC 19330051920289)
#    - add CRUZ / PALMA / VALERIA (NC 19330051920251) in group 4ASV
# ---------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Rescatables")

# Remove the GARCIA/LINARES/ANDRES row (originally row 2); every row
# below shifts up by one.
$wsR.Rows.Item(2).Delete()

# The KARLA MICHEL row was originally row 8; after the delete above it
# is row 7.
$wsR.Rows.Item(7).Delete()

# Insert a fresh row at position 3 for the new student and fill it in.
$wsR.Rows.Item(3).Insert()
$wsR.Range("A3").Value = 19330051920251
$wsR.Range("B3").Value = "CRUZ"
$wsR.Range("C3").Value = "PALMA"
$wsR.Range("D3").Value = "VALERIA"
$wsR.Range("E3").Value = "FÍSICA I"
$wsR.Range("F3").Value = "4ASV"
$wsR.Range("G3").Value = 2

# ---------------------------------------------------------------------
# 2) Estadisticos 1P - row 6 (group 4ALCV) stats refresh
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D6").Value = 0
$ws1.Range("F6").Value = 22
$ws1.Range("G6").Value = 66.67

# ---------------------------------------------------------------------
# 3) Estadisticos 2P - row 6 (group 4ALCV) stats refresh
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 3
$ws2.Range("F6").Value = 30
$ws2.Range("G6").Value = 90.91
$ws2.Range("H6").Value = 7.2

# ---------------------------------------------------------------------
# 4) Estadisticos Final - row 6 (group 4ALCV) stats refresh
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 3
$ws3.Range("F6").Value = 30
$ws3.Range("G6").Value = 90.91
$ws3.Range("H6").Value = 7.5
